# Update countries & provincias Spain
# Applies the 2020-06-15 18:03 COVID data refresh:
#   - Updated case/death counters for a number of countries
#   - Re-sort ripple that swapped two adjacent row pairs
#     (Islas Malvinas <-> Groenlandia, Santa Sede <-> Islas Turcas y Caicos)
#   - Refreshed "Datos actualizados" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Cells.Item(4,2).Value = 2166685
$ws.Cells.Item(4,3).Value = 4457
$ws.Cells.Item(4,4).Value = 870080
$ws.Cells.Item(4,5).Value = 1178703
$ws.Cells.Item(4,7).Value = 44
$ws.Cells.Item(4,8).Value = 117902

# --- Row 7: India ---
$ws.Cells.Item(7,2).Value = 336185
$ws.Cells.Item(7,3).Value = 3402
$ws.Cells.Item(7,4).Value = 172313
$ws.Cells.Item(7,5).Value = 154282
$ws.Cells.Item(7,7).Value = 70
$ws.Cells.Item(7,8).Value = 9590

# --- Row 8: Reino Unido ---
$ws.Cells.Item(8,2).Value = 296857
$ws.Cells.Item(8,3).Value = 968
$ws.Cells.Item(8,7).Value = 38
$ws.Cells.Item(8,8).Value = 41736

# --- Row 13: Alemania ---
$ws.Cells.Item(13,2).Value = 187843
$ws.Cells.Item(13,3).Value = 172
$ws.Cells.Item(13,5).Value = 6368
$ws.Cells.Item(13,7).Value = 5
$ws.Cells.Item(13,8).Value = 8875

# --- Row 20: Canada ---
$ws.Cells.Item(20,2).Value = 99070
$ws.Cells.Item(20,3).Value = 283
$ws.Cells.Item(20,4).Value = 60524
$ws.Cells.Item(20,5).Value = 30372
$ws.Cells.Item(20,7).Value = 28
$ws.Cells.Item(20,8).Value = 8174

# --- Row 33: Singapur ---
$ws.Cells.Item(33,4).Value = 30366
$ws.Cells.Item(33,5).Value = 10426

# --- Row 45: Republica Dominicana ---
$ws.Cells.Item(45,2).Value = 23271
$ws.Cells.Item(45,3).Value = 309
$ws.Cells.Item(45,4).Value = 14025
$ws.Cells.Item(45,5).Value = 8641
$ws.Cells.Item(45,7).Value = 13
$ws.Cells.Item(45,8).Value = 605

# --- Row 56: Kazajistan ---
$ws.Cells.Item(56,4).Value = 9376
$ws.Cells.Item(56,5).Value = 5352
$ws.Cells.Item(56,7).Value = 4
$ws.Cells.Item(56,8).Value = 81

# --- Row 61: Moldavia ---
$ws.Cells.Item(61,2).Value = 11879
$ws.Cells.Item(61,3).Value = 139
$ws.Cells.Item(61,5).Value = 4674
$ws.Cells.Item(61,7).Value = 5
$ws.Cells.Item(61,8).Value = 411

# --- Row 63: Chequia ---
$ws.Cells.Item(63,2).Value = 10044
$ws.Cells.Item(63,3).Value = 20
$ws.Cells.Item(63,4).Value = 7295
$ws.Cells.Item(63,5).Value = 2419

# --- Row 92: Grecia ---
$ws.Cells.Item(92,2).Value = 3134
$ws.Cells.Item(92,3).Value = 13
$ws.Cells.Item(92,5).Value = 1576
$ws.Cells.Item(92,7).Value = 1
$ws.Cells.Item(92,8).Value = 184

# --- Row 103: Sri Lanka ---
$ws.Cells.Item(103,2).Value = 1902
$ws.Cells.Item(103,3).Value = 13
$ws.Cells.Item(103,5).Value = 549

# --- Row 117: Zambia ---
$ws.Cells.Item(117,2).Value = 1382
$ws.Cells.Item(117,3).Value = 24
$ws.Cells.Item(117,4).Value = 1142
$ws.Cells.Item(117,5).Value = 229

# --- Row 126: Republica de Chipre ---
$ws.Cells.Item(126,2).Value = 985
$ws.Cells.Item(126,3).Value = 2
$ws.Cells.Item(126,5).Value = 160

# --- Row 128: Jordania ---
$ws.Cells.Item(128,2).Value = 979
$ws.Cells.Item(128,3).Value = 18
$ws.Cells.Item(128,4).Value = 692
$ws.Cells.Item(128,5).Value = 278

# --- Rows 206-207 swap position: Islas Malvinas <-> Groenlandia ---
# (counters are identical between the two, so only the country names move)
$ws.Cells.Item(206,1).Value = "Groenlandia"
$ws.Cells.Item(207,1).Value = "Islas Malvinas"

# --- Rows 208-209 swap position: Santa Sede <-> Islas Turcas y Caicos ---
$ws.Cells.Item(208,1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(208,4).Value = 11
$ws.Cells.Item(208,8).Value = 1

$ws.Cells.Item(209,1).Value = "Santa Sede"
$ws.Cells.Item(209,4).Value = 12
$ws.Cells.Item(209,8).Value = 0

# --- Refresh the "Datos actualizados" timestamp (shared string not bound to
#     a live cell, so Find/Replace is used to reach it) ---
$ws.Cells.Replace("16:46", "18:03")
